$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '22.419.63'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '1.565.14'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("E5").Value = '  -0.17%  '
$ws.Range("D6").Value = '285.74'
$ws.Range("E6").Value = '  -1.97%  '
$ws.Range("D7").Value = '0.3626'
$ws.Range("E7").Value = '  -2.78%  '
$ws.Range("D8").Value = '48.60'
$ws.Range("E8").Value = '  -2.97%  '
$ws.Range("D9").Value = '0.3341'
$ws.Range("E9").Value = '  -1.42%  '
$ws.Range("D10").Value = '1.125'
$ws.Range("E10").Value = '  -1.33%  '
$ws.Range("D11").Value = '0.07393'
$ws.Range("E11").Value = '  -2.42%  '
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("D13").Value = '20.74'
$ws.Range("E13").Value = '  -2.72%  '
$ws.Range("D14").Value = '5.942'
$ws.Range("E14").Value = '  -0.79%  '
$ws.Range("D15").Value = '6.899'
$ws.Range("E15").Value = '  -0.90%  '
$ws.Range("D16").Value = '1.564.61'
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").Value = '0.00001103'
$ws.Range("E17").Value = '  -1.48%  '
$ws.Range("D18").Value = '88.20'
$ws.Range("E18").Value = '  -3.03%  '
$ws.Range("D19").Value = '0.06692'
$ws.Range("E19").Value = '  -0.63%  '
$ws.Range("D21").Value = '6.342'
$ws.Range("E21").Value = '  +0.92%  '
$ws.Range("D22").Value = '16.16'
$ws.Range("E22").Value = '  -1.04%  '
$ws.Range("E23").Value = '  -0.79%  '
$ws.Range("D24").Value = '22.405.61'
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").Value = '2.403'
$ws.Range("E25").Value = '  +2.87%  '
$ws.Range("D26").Value = '2.563'
$ws.Range("E26").Value = '  -3.69%  '
$ws.Range("D27").Value = '150.48'
$ws.Range("E27").Value = '  +1.40%  '
$ws.Range("D28").Value = '19.35'
$ws.Range("E28").Value = '  -3.66%  '
$ws.Range("D29").Value = '5.011'
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").Value = '123.55'
$ws.Range("E30").Value = '  -1.49%  '
$ws.Range("D31").Value = '1.739.88'
$ws.Range("E31").Value = '  -0.73%  '
$ws.Range("D32").Value = '1.060'
$ws.Range("E32").Value = '  +1.16%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '6.105'
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("B34").Value = 'WEMIXTOKEN'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '1.995'
$ws.Range("E34").Value = '  +1.16%  '
$ws.Range("D35").Value = '9.797'
$ws.Range("E35").Value = '  -0.35%  '
$ws.Range("D36").Value = '0.08272'
$ws.Range("E36").Value = '  -1.30%  '
$ws.Range("D37").Value = '0.02404'
$ws.Range("E37").Value = '  -2.84%  '
$ws.Range("D38").Value = '0.2223'
$ws.Range("E38").Value = '  -2.53%  '
$ws.Range("D39").Value = '0.06389'
$ws.Range("E39").Value = '  -2.02%  '
$ws.Range("D40").Value = '1.297'
$ws.Range("E40").Value = '  -5.77%  '
$ws.Range("D41").Value = '5.330'
$ws.Range("E41").Value = '  -2.41%  '
$ws.Range("D42").Value = '11.14'
$ws.Range("E42").Value = '  -1.10%  '
$ws.Range("D43").Value = '0.6088'
$ws.Range("E43").Value = '  -2.17%  '
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").Value = '13.80'
$ws.Range("E45").Value = '  -1.41%  '
$ws.Range("D46").Value = '3.761'
$ws.Range("E46").Value = '  -1.33%  '
$ws.Range("D47").Value = '0.5795'
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("D48").Value = '2.019'
$ws.Range("E48").Value = '  -2.57%  '
$ws.Range("D49").Value = '124.19'
$ws.Range("E49").Value = '  -4.09%  '
$ws.Range("E50").Value = '  +0.19%  '
$ws.Range("D51").Value = '0.07201'
$ws.Range("E51").Value = '  -1.58%  '
